# Generate Report for Handoff
# -----------------------------------------------------------------------
# This script reproduces a localization-status "handoff" refresh:
#   - the existing in-flight file (9e350708-...md) gets a newer handoff
#     timestamp and a new content hash (so its GUID + xlf hash change)
#   - a brand-new file (ffff47806cfd-...md) shows up "Ready for handoff"
#     and gets appended as row 3 on every sheet
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- new identifiers / values -----------------------------------------
$newMdGuid   = "1761058d-58d9-4c45-a20e-e70b262a33d9"
$newMdGuid2  = "ffff47806cfd-ab1d-46c5-9e26-ae4d4801426c"
$newMdFile   = $newMdGuid + ".md"
$newMdFile2  = $newMdGuid2 + ".md"
$newXlfHash  = "07e8a009d35b56b1d145fb545e9bffebf8cf6370"
$newZhCnXlf  = $newMdGuid + "." + $newXlfHash + ".zh-cn.xlf"
$newDeDeXlf  = $newMdGuid + "." + $newXlfHash + ".de-de.xlf"

$overviewDate = "2016-03-22 06:53:26"
$zhCnDate     = "2016-03-22 06:53:18"
$deDeDate     = "2016-03-22 06:53:26"
$epochDate    = "0001-01-01 00:00:00"
$readyStatus  = "Ready for handoff"
$includeVal   = "Include"
$dotMd        = ".md"

$mdUrlBase     = "https://github.com/OpenLocalizationTest/oltest/blob/c276526862dee5d1402ed381ed742181d1ddfaac/e2e/"
$zhCnXlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f31d1009de14890d31d50978a2ac89ff0a3575be/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/"
$deDeXlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b7666ebb039b18ba3bd85c2d2f280b6d5ff4cc67/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/"

$dateFormat = "yyyy-mm-dd HH:mm:ss"
$hlColor = 6591981     # cornflower blue FF6495ED, expressed RGB(100,149,237)
$hlColorBgr = 15570276 # same colour but through the BGR-ish channel this host expects

function Set-CellHyperlink($ws, $cellRef, $address, $display) {
    $rng = $ws.Range($cellRef)
    $rng.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($rng, $address, "", "", $display)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = 2
    $rng.Font.Color = $hlColorBgr
}

function Set-CellDate($ws, $cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $value
    $rng.NumberFormat = $dateFormat
}

# =========================================================================
# Sheet "Overview"
# =========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2: refresh the in-flight file's link + handoff date
Set-CellHyperlink $wsOverview "A2" ($mdUrlBase + $newMdFile) $newMdFile
Set-CellDate $wsOverview "D2" $overviewDate

# Row 3: the newly-ready file
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus
Set-CellHyperlink $wsOverview "A3" ($mdUrlBase + $newMdFile2) $newMdFile2
Set-CellDate $wsOverview "D3" $overviewDate

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2: refresh links + handoff datetime
Set-CellHyperlink $wsZhCn "A2" ($mdUrlBase + $newMdFile) $newMdFile
Set-CellHyperlink $wsZhCn "D2" ($zhCnXlfUrlBase + $newZhCnXlf) $newZhCnXlf
Set-CellDate $wsZhCn "E2" $zhCnDate

# Row 3: the newly-ready file (same shape as row 2)
$wsZhCn.Range("B3").Value = $dotMd
$wsZhCn.Range("C3").Value = $readyStatus
Set-CellHyperlink $wsZhCn "A3" ($mdUrlBase + $newMdFile2) $newMdFile2
Set-CellHyperlink $wsZhCn "D3" ($zhCnXlfUrlBase + $newZhCnXlf) $newZhCnXlf
Set-CellDate $wsZhCn "E3" $zhCnDate
Set-CellDate $wsZhCn "H3" $epochDate
$wsZhCn.Range("J3").Value = $includeVal

# =========================================================================
# Sheet "de-de"
# =========================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2: refresh links + handoff datetime
Set-CellHyperlink $wsDeDe "A2" ($mdUrlBase + $newMdFile) $newMdFile
Set-CellHyperlink $wsDeDe "D2" ($deDeXlfUrlBase + $newDeDeXlf) $newDeDeXlf
Set-CellDate $wsDeDe "E2" $deDeDate

# Row 3: the newly-ready file (same shape as row 2)
$wsDeDe.Range("B3").Value = $dotMd
$wsDeDe.Range("C3").Value = $readyStatus
Set-CellHyperlink $wsDeDe "A3" ($mdUrlBase + $newMdFile2) $newMdFile2
Set-CellHyperlink $wsDeDe "D3" ($deDeXlfUrlBase + $newDeDeXlf) $newDeDeXlf
Set-CellDate $wsDeDe "E3" $deDeDate
Set-CellDate $wsDeDe "H3" $epochDate
$wsDeDe.Range("J3").Value = $includeVal

Write-Host "Handoff report regenerated."
